$d = $word.ActiveDocument

# Locate the anchor paragraph (the last paragraph of section 2.1, ending with
# "...response ratio keeps increasing.") so the new "2.2 Round Robin" content
# is inserted right after it, regardless of exact paragraph index.
$anchorText = "中长作业可能会饿死的问题，因为随着等待时间的增长，响应比也会越来越高。"
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Contains($anchorText)) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "anchor paragraph not found"
}

function Insert-ParaXml($afterIndex, $xml) {
    $p = $d.Paragraphs.Item($afterIndex)
    $r = $p.Range.Duplicate
    $r.Collapse(0)
    [void]$r.InsertParagraphAfter()
    $newp = $d.Paragraphs.Item($afterIndex + 1)
    [void]$newp.Range.InsertXML($xml)
}

$idx = $anchorIndex

$xmlEmpty = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>
'@
Insert-ParaXml $idx $xmlEmpty
$idx = $idx + 1

$xmlHeading = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">2.2 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>时间片轮转</w:t></w:r></w:p>
'@
Insert-ParaXml $idx $xmlHeading
$idx = $idx + 1

$xmlBody = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>将所有就绪进程按</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> FCFS </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>的原则排成一个队列，每次调度时，把</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> CPU </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>分配给队首进程，该进程可以执行一个时间片。当时间片用完时，由计时器发出时钟中断，调度程序便停止该进程的执行，并将它送往就绪队列的末尾，同时继续把</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> CPU </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>分配给队首的进程。</w:t></w:r></w:p>
'@
Insert-ParaXml $idx $xmlBody
$idx = $idx + 1

$xmlLast = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>时间片轮转算法的效率和时间片的大小有很大关系。因为每次进程切换都要保存进程的信息并且载入新进程的信息，如果时间片太小，进程切换太频繁，在进程切换上就会花过多时间。</w:t></w:r></w:p>
'@
Insert-ParaXml $idx $xmlLast
